$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Row 72: update the "trying moogfilter again" measurement
$ws.Range("B72").Value = 20608

# Row 73: new entry "sat: no slope option at all (24db/oct forced)"
$ws.Range("A73").Value = "sat: no slope option at all (24db/oct forced)"
$ws.Range("B73").Value = 20380
$ws.Range("D73").Value = 21276
$ws.Range("F73").Value = "ok now we're talking. Moog is back baby"

# New trailing rows 74-77 carrying forward the baseline "D" value
$ws.Range("D74").Value = 21276
$ws.Range("D75").Value = 21276
$ws.Range("D76").Value = 21276
$ws.Range("D77").Value = 21276

# Update the selection to match the edited cell
$ws.Range("A73").Select()
